$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(584).Insert()
$ws.Range("A584").Value = 10
$ws.Range("B584").Value = "Vega Modelo de Temuco"
$ws.Range("C584").Value = "La Araucanía"
$ws.Range("D584").Value = 45212
$ws.Range("E584").Value = 9
$ws.Range("F584").Value = "Fruta"
$ws.Range("G584").Value = 100108
$ws.Range("H584").Value = "Tropicales y subtropicales"
$ws.Range("I584").Value = 100108002
$ws.Range("J584").Value = "Mango"
$ws.Range("K584").Value = "Sin especificar"
$ws.Range("L584").Value = "Primera"
$ws.Range("M584").Value = 500
$ws.Range("N584").Value = 10000
$ws.Range("O584").Value = 10000
$ws.Range("P584").Value = 10000
$ws.Range("Q584").Value = '$/bandeja 4 kilos'
$ws.Range("R584").Value = "Brasil"
$ws.Range("S584").Value = 2500
$ws.Range("T584").Value = 4
